# Regenerate the 15 lattice-multiplication exercise cells (5 rows x 3 cols)
# with a new set of multiplication problems, per the target revision.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New "A B x C D" pairs, in row-major order (row1: col1,col2,col3; row2: ...)
$pairs = @(
    "36 x 58", "87 x 94", "86 x 33",
    "57 x 75", "26 x 14", "52 x 18",
    "49 x 83", "38 x 27", "51 x 39",
    "54 x 59", "43 x 90", "14 x 49",
    "49 x 21", "90 x 41", "57 x 84"
)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$idx = 0
for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $pair = $pairs[$idx]
        $idx = $idx + 1

        $nums = $pair -split " x "
        $a = $nums[0]
        $b = $nums[1]

        $breakdown = "  {0}    {1}" -f $b.Substring(0,1), $b.Substring(1,1)
        $row1 = "{0}|    |" -f $a.Substring(0,1)
        $row2 = "{0}|    |" -f $a.Substring(1,1)

        $xml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + `
               '<w:t>' + $pair + '</w:t><w:br/>' + `
               '<w:t xml:space="preserve">' + $breakdown + '</w:t><w:br/>' + `
               '<w:t xml:space="preserve">  ----</w:t><w:br/>' + `
               '<w:t>' + $row1 + '</w:t><w:br/>' + `
               '<w:t>' + $row2 + '</w:t>' + `
               '</w:r></w:p>'

        $cell = $t.Cell($r, $c)
        [void]$cell.Range.InsertXML($xml)
    }
}
